$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[Shuang-Jiang%Zhou%NULL%0,         Li-Gang%Zhang%NULL%0,         Lei-Lei%Wang%NULL%0,         Zhao-Chang%Guo%NULL%0,         Jing-Qi%Wang%NULL%0,         Jin-Cheng%Chen%NULL%0,         Mei%Liu%NULL%0,         Xi%Chen%NULL%0,         Jing-Xu%Chen%chenjx1110@163.com%0]"
$ws.Range("E3").Value = "[Benjamin%Oosterhoff%Benjamin.oosterhoff@montana.edu%0,         Cara A.%Palmer%NULL%0,         Jenna%Wilson%NULL%0,         Natalie%Shook%NULL%0]"
$ws.Range("E4").Value = "[İsmail%Seçer%ismailsecer84@gmail.com%0,         Sümeyye%Ulaş%NULL%0,         Sümeyye%Ulaş%NULL%0]"
$ws.Range("E5").Value = "[Miao%Qu%NULL%1,         Kun%Yang%NULL%0,         Yujia%Cao%NULL%1,         Mei Hong%Xiu%xiumeihong97@163.com%1,         Xiang Yang%Zhang%zhangxy@psych.ac.cn%2,         Xiang Yang%Zhang%zhangxy@psych.ac.cn%0]"

Write-Output "Updated E2:E5 with new author reference strings (crossref update)"
